# TUM Lite Enchanted 1.3 - Rebalance Creatures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Row 6: creature cost tweak (column I)
$ws.Range("I6").Value = 17030

# Row 17: creature cost tweaks (columns I and J)
$ws.Range("I17").Value = 25790
$ws.Range("J17").Value = 23900

# Row 28: recalibrated breakpoint formulas (columns C, I, J)
$ws.Range("C28").Formula = "=(150*20)+(200*16)+(290*13)+(500*8)+(550*6)+(2000*4)+(10000*2)"
$ws.Range("I28").Formula = "=(85*38)+(200*22)+(310*18)+(520*8)+(755*6)+(2000*4)+(7400*2)"
$ws.Range("J28").Formula = "=(110*28)+(230*20)+(330*16)+(450*8)+(890*6)+(1860*4)+(8000*2)"
